$wb = $excel.ActiveWorkbook

# --- Update Hoja1!A1 text block with new conversion rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.94 = 10829.41 pesos`n✅ 10829.41 pesos = 2.92 = 951.1 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update tasas sheet rates ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 340
$ws2.Range("O10").Value = 3682
$ws2.Range("N12").Value = 3711.91
$ws2.Range("O12").Value = 326
